$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 (RF) - update B3:J3
$ws.Range("B3").Value = 0.788
$ws.Range("C3").Value = 0.883
$ws.Range("D3").Value = 0.636
$ws.Range("E3").Value = 0.898
$ws.Range("F3").Value = 0.89
$ws.Range("G3").Value = 0.112
$ws.Range("H3").Value = 0.335
$ws.Range("I3").Value = 0.241
$ws.Range("J3").Value = 0.972

# Row 4 (NN) - update E4:J4 only
$ws.Range("E4").Value = 0.803
$ws.Range("F4").Value = 0.787
$ws.Range("G4").Value = 0.216
$ws.Range("H4").Value = 0.465
$ws.Range("I4").Value = 0.349
$ws.Range("J4").Value = 0.936

# Row 5 (RNN) - update E5:J5 only
$ws.Range("E5").Value = 0.613
$ws.Range("F5").Value = 0.596
$ws.Range("G5").Value = 0.422
$ws.Range("H5").Value = 0.65
$ws.Range("I5").Value = 0.468
$ws.Range("J5").Value = 0.784
